# Build site at 2022-09-26 16:07:08 UTC
#
# The "Docentes responsáveis:" section used to have its value on its own
# row (row 13, with no label in column A) directly below the label row
# (row 12). That value row is removed entirely (the whole row is deleted,
# shifting every following row up by one), and the professor's name is
# instead written into the "Objetivos:" row. A few other rows further
# down then end up carrying values that belong to neighbouring rows
# (matching the published sheet exactly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the stray value-only row (old row 13) - shifts rows 14-25 up by one.
$ws.Rows(13).Delete()

# 2) After the shift, patch the cells whose text content changed.
$ws.Range("B10:C10").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15:C15").Value = "01/01/2018"
$ws.Range("B18:C18").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("B19:C19").Value = "Por meio de aulas presenciais, com apresentação dos fundamentos e exemplos ou casos, e também pela apresentação de trabalhos em equipes.Justificativa: adequação do método de avaliação."
$ws.Range("B20:C20").Value = "A Avaliação será: MF = (P1 + P2)/2; Onde: P1: Trabalho; P2: Trabalho. Poderá haver também prova individual sobre os fundamentos."
$ws.Range("B21:C21").Value = "Prova de exame."
